# Updated cryptos list data (price / 1h volume change) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.582.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.865.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.98"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.03"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.863.58"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.47"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000266"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.05"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.514.28"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.853.47"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.743.65"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.44%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "470.68"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.733"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.54%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.49"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.22"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.11%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.017.64"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.67"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.50"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.30"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.33"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.829.92"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.71"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.94%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.140"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.02"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.89%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.314"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.71"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "417.06"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.38%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.11"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.13%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000291"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.70%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.23"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.59%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.13"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.01%  "
